$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("Z2").Font.Size = 11
$ws1.Range("Z2").NumberFormat = "0.00"
